$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 5).Value = 0.01547678482276585
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 0.01092523045407989
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 0.003231017770597738
$ws.Cells.Item(3, 5).Value = 0.01647528706939591
$ws.Cells.Item(3, 6).Value = 0.004116766467065868
$ws.Cells.Item(3, 7).Value = 0.03448275862068965
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 0.008885298869143781
$ws.Cells.Item(4, 2).Value = 0.007343941248470012
$ws.Cells.Item(4, 3).Value = 0.02445265851578049
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = 0.0144782825761358
$ws.Cells.Item(4, 6).Value = 0.02357784431137723
$ws.Cells.Item(4, 7).Value = 0.05565039262546945
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 0.002423263327948304
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(5, 2).Value = 0.2288861689106495
$ws.Cells.Item(5, 3).Value = 0.1595109468296855
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(5, 5).Value = 0.005491762356465302
$ws.Cells.Item(5, 6).Value = 0.06025449101796427
$ws.Cells.Item(5, 7).Value = 0.02799590303857972
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 0.0008077544426494346
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(6, 5).Value = 0.006490264603095357
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 0.0102424035506999
$ws.Cells.Item(6, 9).Value = 0
$ws.Cells.Item(6, 10).Value = 0.001615508885298869
$ws.Cells.Item(7, 5).Value = 0.01248127808287569
$ws.Cells.Item(7, 8).Value = 0.007874015748031496
$ws.Cells.Item(7, 10).Value = 0.01130856219709209
$ws.Cells.Item(7, 11).Value = 0.006097560975609756
$ws.Cells.Item(8, 2).Value = 0.01101591187270501
$ws.Cells.Item(8, 3).Value = 0.08984930338356603
$ws.Cells.Item(8, 4).Value = 0.01612903225806452
$ws.Cells.Item(8, 5).Value = 0.04493260109835256
$ws.Cells.Item(8, 6).Value = 0.2889221556886206
$ws.Cells.Item(8, 7).Value = 0.3287811539774632
$ws.Cells.Item(8, 8).Value = 0.01837270341207349
$ws.Cells.Item(8, 9).Value = 0
$ws.Cells.Item(8, 10).Value = 0.02423263327948305
$ws.Cells.Item(8, 11).Value = 0.02134146341463415
$ws.Cells.Item(9, 5).Value = 0.01647528706939591
$ws.Cells.Item(9, 10).Value = 0.004846526655896607
$ws.Cells.Item(9, 11).Value = 0.001524390243902439
$ws.Cells.Item(10, 2).Value = 0
$ws.Cells.Item(10, 3).Value = 0.003127665624111458
$ws.Cells.Item(10, 4).Value = 0.03629032258064514
$ws.Cells.Item(10, 5).Value = 0.1387918122815784
$ws.Cells.Item(10, 6).Value = 0.01908682634730538
$ws.Cells.Item(10, 7).Value = 0.0307272106520997
$ws.Cells.Item(10, 8).Value = 0.02362204724409449
$ws.Cells.Item(10, 9).Value = 0.005673758865248227
$ws.Cells.Item(10, 10).Value = 0.04281098546042001
$ws.Cells.Item(10, 11).Value = 0.04115853658536587
$ws.Cells.Item(12, 2).Value = 0
$ws.Cells.Item(12, 4).Value = 0.008064516129032258
$ws.Cells.Item(12, 8).Value = 0.03412073490813648
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 10).Value = 0.02180936995153475
$ws.Cells.Item(12, 11).Value = 0.04115853658536587
$ws.Cells.Item(13, 3).Value = 0.01393232868922379
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0.000998502246630055
$ws.Cells.Item(13, 6).Value = 0.01085329341317365
$ws.Cells.Item(13, 7).Value = 0.002731307613519973
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(14, 2).Value = 0.2564259485924119
$ws.Cells.Item(14, 3).Value = 0.04492465169178286
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = 0.008986520219670495
$ws.Cells.Item(14, 6).Value = 0.007110778443113775
$ws.Cells.Item(14, 7).Value = 0.009559576647319904
$ws.Cells.Item(14, 9).Value = 0
$ws.Cells.Item(14, 10).Value = 0.001615508885298869
$ws.Cells.Item(14, 11).Value = 0
$ws.Cells.Item(15, 5).Value = 0.02845731402895657
$ws.Cells.Item(16, 2).Value = 0
$ws.Cells.Item(16, 3).Value = 0.03724765425078194
$ws.Cells.Item(16, 4).Value = 0.006048387096774193
$ws.Cells.Item(16, 5).Value = 0.03544682975536698
$ws.Cells.Item(16, 6).Value = 0.03143712574850296
$ws.Cells.Item(16, 7).Value = 0.03721406623420963
$ws.Cells.Item(16, 8).Value = 0.01312335958005249
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 10).Value = 0.01696284329563813
$ws.Cells.Item(16, 11).Value = 0.009146341463414634
$ws.Cells.Item(17, 5).Value = 0.000998502246630055
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(17, 7).Value = 0.002048480710139979
$ws.Cells.Item(17, 10).Value = 0.0008077544426494346
$ws.Cells.Item(18, 2).Value = 0
$ws.Cells.Item(18, 4).Value = 0.002016129032258064
$ws.Cells.Item(18, 8).Value = 0.005249343832020997
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 10).Value = 0.001615508885298869
$ws.Cells.Item(18, 11).Value = 0.003048780487804878
$ws.Cells.Item(19, 2).Value = 0.01223990208078335
$ws.Cells.Item(19, 3).Value = 0.02018765993744669
$ws.Cells.Item(19, 4).Value = 0
$ws.Cells.Item(19, 5).Value = 0.00599101347978033
$ws.Cells.Item(19, 6).Value = 0.01085329341317365
$ws.Cells.Item(19, 7).Value = 0.01946056674632981
$ws.Cells.Item(19, 9).Value = 0
$ws.Cells.Item(19, 10).Value = 0.0008077544426494346
$ws.Cells.Item(19, 11).Value = 0
$ws.Cells.Item(20, 5).Value = 0.006490264603095357
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(20, 7).Value = 0.009559576647319904
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(21, 3).Value = 0.009667330110889967
$ws.Cells.Item(21, 4).Value = 0
$ws.Cells.Item(21, 5).Value = 0.004493260109835247
$ws.Cells.Item(21, 6).Value = 0.007110778443113775
$ws.Cells.Item(21, 7).Value = 0.006145442130419939
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 10).Value = 0.001615508885298869
$ws.Cells.Item(21, 11).Value = 0
$ws.Cells.Item(22, 4).Value = 0.01209677419354839
$ws.Cells.Item(22, 5).Value = 0.000998502246630055
$ws.Cells.Item(22, 8).Value = 0.02362204724409449
$ws.Cells.Item(22, 9).Value = 0.005673758865248227
$ws.Cells.Item(22, 10).Value = 0.04604200323101774
$ws.Cells.Item(22, 11).Value = 0.05792682926829271
$ws.Cells.Item(23, 5).Value = 0.03145282076884673
$ws.Cells.Item(23, 6).Value = 0.004865269461077844
$ws.Cells.Item(23, 7).Value = 0.06145442130419939
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 10).Value = 0.02584814216478192
$ws.Cells.Item(23, 11).Value = 0.003048780487804878
$ws.Cells.Item(24, 2).Value = 0
$ws.Cells.Item(24, 4).Value = 0.002016129032258064
$ws.Cells.Item(24, 5).Value = 0.0004992511233150275
$ws.Cells.Item(24, 8).Value = 0.002624671916010499
$ws.Cells.Item(24, 9).Value = 0
$ws.Cells.Item(24, 10).Value = 0.004846526655896607
$ws.Cells.Item(24, 11).Value = 0.006097560975609756
$ws.Cells.Item(32, 2).Value = 0.08384332925336616
$ws.Cells.Item(32, 3).Value = 0.03212965595678133
$ws.Cells.Item(32, 4).Value = 0.008064516129032258
$ws.Cells.Item(32, 5).Value = 0.08637044433350013
$ws.Cells.Item(32, 6).Value = 0.01085329341317365
$ws.Cells.Item(32, 7).Value = 0.02594742232843974
$ws.Cells.Item(32, 8).Value = 0.03674540682414698
$ws.Cells.Item(32, 9).Value = 0
$ws.Cells.Item(32, 10).Value = 0.08077544426494346
$ws.Cells.Item(32, 11).Value = 0.04420731707317075
$ws.Cells.Item(33, 2).Value = 0
$ws.Cells.Item(33, 3).Value = 0
$ws.Cells.Item(33, 4).Value = 0.06653225806451608
$ws.Cells.Item(33, 5).Value = 0.1073389915127314
$ws.Cells.Item(33, 7).Value = 0.002389894161829976
$ws.Cells.Item(33, 8).Value = 0.141732283464567
$ws.Cells.Item(33, 9).Value = 0
$ws.Cells.Item(33, 10).Value = 0.1623586429725368
$ws.Cells.Item(33, 11).Value = 0.1600609756097562
$ws.Cells.Item(34, 5).Value = 0.0004992511233150275
$ws.Cells.Item(34, 6).Value = 0.0007485029940119761
$ws.Cells.Item(34, 7).Value = 0.00102424035506999
$ws.Cells.Item(34, 10).Value = 0